# Applies the two changes captured by the commit's OOXML diff:
#
#  1. Slide 16's table (the PLENARY fill-in-the-gaps table) switches its
#     table style from {FA85B958-CD31-4CC4-946B-56E327F7B0FB} to
#     {7774F01B-B663-4134-8224-6DC35AA41125}.
#
#  2. The presentation's theme (ppt/theme/theme1.xml, used by the slide
#     master) swaps its colour scheme from the "Integral" palette to the
#     stock "Office Theme" palette (the two themes already share the same
#     font scheme / format scheme, so only the 12 colour-scheme slots -
#     dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - actually differ).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table: find it by its current (old) style GUID so the
#    edit is robust even if shape/slide ordering shifts.
# ---------------------------------------------------------------------
$oldStyle = "{FA85B958-CD31-4CC4-946B-56E327F7B0FB}"
$newStyle = "{7774F01B-B663-4134-8224-6DC35AA41125}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colour scheme (Integral -> Office) on the slide
#    master's theme. Order of ThemeColorScheme.Item(n) is fixed:
#    1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
#    RGB is the usual COM BGR-packed integer.
# ---------------------------------------------------------------------
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($ci = 1; $ci -le $colorScheme.Count; $ci++) {
    $colorScheme.Item($ci).RGB = $officeColors[$ci - 1]
}
